# "Cleaned up folder deletion" - remove the now-unused "Collected" column
# from the CapitalCommitment sheet. The column sat between "Committed
# Amount" (C) and "SPV" (D), so deleting it shifts SPV/Phone one column
# to the left (D/E) and drops the "Collected" / "\u20b90.00" strings that
# are no longer referenced anywhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the column first (mirrors the real user action of right-clicking
# the column header) and then delete it, shifting everything after it left.
$ws.Columns("D").Select() | Out-Null
$ws.Columns("D").Delete()
